$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "1047483918"
$ws.Range("D16").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E16").Value = "1909"
$ws.Range("F16").Value = 25396

$ws.Range("C17").Value = "1047483918"
$ws.Range("D17").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E17").Value = "1908"
$ws.Range("F17").Value = 33125

$ws.Range("C18").Value = "1047483918"
$ws.Range("D18").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E18").Value = "1907"
$ws.Range("F18").Value = 33125

$ws.Range("C19").Value = "1047483918"
$ws.Range("D19").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E19").Value = "1906"
$ws.Range("F19").Value = 33125

$ws.Range("C20").Value = "1047483918"
$ws.Range("D20").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E20").Value = "1905"
$ws.Range("F20").Value = 33125

$ws.Range("C21").Value = "1047483918"
$ws.Range("D21").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E21").Value = "1904"
$ws.Range("F21").Value = 33125

$ws.Range("C22").Value = "1047483918"
$ws.Range("D22").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E22").Value = "1903"
$ws.Range("F22").Value = 33125

$ws.Range("C23").Value = "1047483918"
$ws.Range("D23").Value = "LUIS ANTONIO RUIZ MUÑOZ"
$ws.Range("E23").Value = "1902"
$ws.Range("F23").Value = 8833

$ws.Range("C24").Value = "1143325267"
$ws.Range("D24").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E24").Value = "1909"
$ws.Range("F24").Value = 25396

$ws.Range("C25").Value = "1143325267"
$ws.Range("D25").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E25").Value = "1908"
$ws.Range("F25").Value = 33125

$ws.Range("C26").Value = "1143325267"
$ws.Range("D26").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E26").Value = "1907"
$ws.Range("F26").Value = 33125

$ws.Range("C27").Value = "1143325267"
$ws.Range("D27").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E27").Value = "1906"
$ws.Range("F27").Value = 33125

$ws.Range("C28").Value = "1143325267"
$ws.Range("D28").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E28").Value = "1905"
$ws.Range("F28").Value = 33125

$ws.Range("C29").Value = "1143325267"
$ws.Range("D29").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E29").Value = "1904"
$ws.Range("F29").Value = 33125

$ws.Range("C30").Value = "1143325267"
$ws.Range("D30").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E30").Value = "1903"
$ws.Range("F30").Value = 33125

$ws.Range("C31").Value = "1143325267"
$ws.Range("D31").Value = "WALTER RODRIGUEZ ROMERIN"
$ws.Range("E31").Value = "1902"
$ws.Range("F31").Value = 8833

$ws.Range("C32").Value = "73000153"
$ws.Range("D32").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E32").Value = "1909"
$ws.Range("F32").Value = 25396

$ws.Range("C33").Value = "73000153"
$ws.Range("D33").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E33").Value = "1908"
$ws.Range("F33").Value = 33125

$ws.Range("C34").Value = "73000153"
$ws.Range("D34").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E34").Value = "1907"
$ws.Range("F34").Value = 33125

$ws.Range("C35").Value = "73000153"
$ws.Range("D35").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E35").Value = "1906"
$ws.Range("F35").Value = 33125

$ws.Range("C36").Value = "73000153"
$ws.Range("D36").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E36").Value = "1905"
$ws.Range("F36").Value = 33125

$ws.Range("C37").Value = "73000153"
$ws.Range("D37").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E37").Value = "1904"
$ws.Range("F37").Value = 33125

$ws.Range("C38").Value = "73000153"
$ws.Range("D38").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E38").Value = "1903"
$ws.Range("F38").Value = 33125

$ws.Range("C39").Value = "73000153"
$ws.Range("D39").Value = "ALEJANDRO ANTONIO CASSIANI AYALA"
$ws.Range("E39").Value = "1902"
$ws.Range("F39").Value = 8833
